$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price column + Volume(1h) label column)
# Leading apostrophe forces text storage so numeric-looking strings
# (e.g. "235.00") aren't silently converted to numbers, matching the
# original inline-string text cells.

$ws.Range("D2").Value = "'235.00"
$ws.Range("D3").Value = "'22.48"
$ws.Range("D4").Value = "'5.300"
$ws.Range("D5").Value = "'0.05613"
$ws.Range("D6").Value = "'3.379"
$ws.Range("D7").Value = "'6.490"
$ws.Range("D8").Value = "'1.069"
$ws.Range("D9").Value = "'0.7820"
$ws.Range("D10").Value = "'0.1398"
$ws.Range("D11").Value = "'0.07396"
$ws.Range("D12").Value = "'0.03152"
$ws.Range("D13").Value = "'0.02988"
$ws.Range("D14").Value = "'0.09261"
$ws.Range("D15").Value = "'0.001655"
$ws.Range("D16").Value = "'3.253"
$ws.Range("D17").Value = "'0.04732"
$ws.Range("E18").Value = "'17OneONE"
$ws.Range("D19").Value = "'0.006202"
$ws.Range("D20").Value = "'0.005222"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.962"
$ws.Range("D27").Value = "'0.0004990"
$ws.Range("E27").Value = "'26UpBotsUBXT"
$ws.Range("D40").Value = "'0.04054"
$ws.Range("D41").Value = "'0.006962"
$ws.Range("E41").Value = "'40KickTokenKICKBestin24h"
$ws.Range("D42").Value = "'0.003500"
$ws.Range("D44").Value = "'0.009429"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.6752"
$ws.Range("D48").Value = "'0.04042"
$ws.Range("E48").Value = "'47BOLOBOLOWorstin24h"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.01010"
